# Re-apply the Normal style to the pre-existing data cells (this is what
# flips applyFont on the shared cell format without introducing a new font).
# Must run before adding the new cells below, since those stay on the
# original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:F1").Style = "Normal"
$ws.Range("A2:E2").Style = "Normal"
$ws.Range("A3:D3").Style = "Normal"

# Add the two new cells containing the "&" character that the commit is about.
$ws.Range("F2").Value = "&Teste"
$ws.Range("G2").Value = "&Outro &Teste"

# Give the used columns an explicit width, matching the saved file as
# closely as this COM layer's rounding allows.
$ws.Columns.Item(1).ColumnWidth = 7.5

# Move the active selection to G3, matching the saved file's cursor position.
$ws.Range("G3").Select()
